$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, matching style of E1 (bold, bordered, centered header style)
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "time_taken"

# time_taken values for data rows 2..287
$timeTaken = @(
    "2021-10-05 13:38:40.988666",
    "2021-10-05 13:38:40.988677",
    "2021-10-05 13:38:40.988680",
    "2021-10-05 13:38:40.988683",
    "2021-10-05 13:38:40.988686",
    "2021-10-05 13:38:40.988689",
    "2021-10-05 13:38:40.988691",
    "2021-10-05 13:38:40.988694",
    "2021-10-05 13:38:40.988697",
    "2021-10-05 13:38:40.988699",
    "2021-10-05 13:38:40.988702",
    "2021-10-05 13:38:40.988704",
    "2021-10-05 13:38:40.988707",
    "2021-10-05 13:38:40.988709",
    "2021-10-05 13:38:40.988712",
    "2021-10-05 13:38:40.988714",
    "2021-10-05 13:38:40.988717",
    "2021-10-05 13:38:40.988719",
    "2021-10-05 13:38:40.988722",
    "2021-10-05 13:38:40.988724",
    "2021-10-05 13:38:40.988727",
    "2021-10-05 13:38:40.988729",
    "2021-10-05 13:38:40.988731",
    "2021-10-05 13:38:40.988734",
    "2021-10-05 13:38:40.988737",
    "2021-10-05 13:38:40.988739",
    "2021-10-05 13:38:40.988742",
    "2021-10-05 13:38:40.988744",
    "2021-10-05 13:38:40.988747",
    "2021-10-05 13:38:40.988749",
    "2021-10-05 13:38:40.988752",
    "2021-10-05 13:38:40.988754",
    "2021-10-05 13:38:40.988757",
    "2021-10-05 13:38:40.988759",
    "2021-10-05 13:38:40.988762",
    "2021-10-05 13:38:40.988764",
    "2021-10-05 13:38:40.988767",
    "2021-10-05 13:38:40.988769",
    "2021-10-05 13:38:40.988772",
    "2021-10-05 13:38:40.988774",
    "2021-10-05 13:38:40.988777",
    "2021-10-05 13:38:40.988779",
    "2021-10-05 13:38:40.988782",
    "2021-10-05 13:38:40.988784",
    "2021-10-05 13:38:40.988787",
    "2021-10-05 13:38:40.988789",
    "2021-10-05 13:38:40.988792",
    "2021-10-05 13:38:40.988794",
    "2021-10-05 13:38:40.988796",
    "2021-10-05 13:38:40.988799",
    "2021-10-05 13:38:40.988801",
    "2021-10-05 13:38:40.988804",
    "2021-10-05 13:38:40.988807",
    "2021-10-05 13:38:40.988809",
    "2021-10-05 13:38:40.988812",
    "2021-10-05 13:38:40.988814",
    "2021-10-05 13:38:40.988817",
    "2021-10-05 13:38:40.988819",
    "2021-10-05 13:38:40.988821",
    "2021-10-05 13:38:40.988824",
    "2021-10-05 13:38:40.988826",
    "2021-10-05 13:38:40.988829",
    "2021-10-05 13:38:40.988831",
    "2021-10-05 13:38:40.988833",
    "2021-10-05 13:38:40.988837",
    "2021-10-05 13:38:40.988840",
    "2021-10-05 13:38:40.988842",
    "2021-10-05 13:38:40.988845",
    "2021-10-05 13:38:40.988847",
    "2021-10-05 13:38:40.988850",
    "2021-10-05 13:38:40.988852",
    "2021-10-05 13:38:40.988855",
    "2021-10-05 13:38:40.988857",
    "2021-10-05 13:38:40.988859",
    "2021-10-05 13:38:40.988862",
    "2021-10-05 13:38:40.988864",
    "2021-10-05 13:38:40.988869",
    "2021-10-05 13:38:40.988871",
    "2021-10-05 13:38:40.988874",
    "2021-10-05 13:38:40.988877",
    "2021-10-05 13:38:40.988879",
    "2021-10-05 13:38:40.988881",
    "2021-10-05 13:38:40.988884",
    "2021-10-05 13:38:40.988886",
    "2021-10-05 13:38:40.988889",
    "2021-10-05 13:38:40.988891",
    "2021-10-05 13:38:40.988894",
    "2021-10-05 13:38:40.988896",
    "2021-10-05 13:38:40.988899",
    "2021-10-05 13:38:40.988901",
    "2021-10-05 13:38:40.988904",
    "2021-10-05 13:38:40.988906",
    "2021-10-05 13:38:40.988910",
    "2021-10-05 13:38:40.988913",
    "2021-10-05 13:38:40.988915",
    "2021-10-05 13:38:40.988918",
    "2021-10-05 13:38:40.988920",
    "2021-10-05 13:38:40.988923",
    "2021-10-05 13:38:40.988925",
    "2021-10-05 13:38:40.988928",
    "2021-10-05 13:38:40.988930",
    "2021-10-05 13:38:40.988933",
    "2021-10-05 13:38:40.988935",
    "2021-10-05 13:38:40.988938",
    "2021-10-05 13:38:40.988940",
    "2021-10-05 13:38:40.988943",
    "2021-10-05 13:38:40.988945",
    "2021-10-05 13:38:40.988948",
    "2021-10-05 13:38:40.988953",
    "2021-10-05 13:38:40.988956",
    "2021-10-05 13:38:40.988959",
    "2021-10-05 13:38:40.988961",
    "2021-10-05 13:38:40.988964",
    "2021-10-05 13:38:40.988966",
    "2021-10-05 13:38:40.988969",
    "2021-10-05 13:38:40.988971",
    "2021-10-05 13:38:40.988974",
    "2021-10-05 13:38:40.988976",
    "2021-10-05 13:38:40.988979",
    "2021-10-05 13:38:40.988981",
    "2021-10-05 13:38:40.988983",
    "2021-10-05 13:38:40.988986",
    "2021-10-05 13:38:40.988988",
    "2021-10-05 13:38:40.988991",
    "2021-10-05 13:38:40.988993",
    "2021-10-05 13:38:40.988995",
    "2021-10-05 13:38:40.988998",
    "2021-10-05 13:38:40.989000",
    "2021-10-05 13:38:40.989005",
    "2021-10-05 13:38:40.989008",
    "2021-10-05 13:38:40.989011",
    "2021-10-05 13:38:40.989013",
    "2021-10-05 13:38:40.989016",
    "2021-10-05 13:38:40.989018",
    "2021-10-05 13:38:40.989020",
    "2021-10-05 13:38:40.989023",
    "2021-10-05 13:38:40.989025",
    "2021-10-05 13:38:40.989028",
    "2021-10-05 13:38:40.989031",
    "2021-10-05 13:38:40.989033",
    "2021-10-05 13:38:40.989035",
    "2021-10-05 13:38:40.989038",
    "2021-10-05 13:38:40.989040",
    "2021-10-05 13:38:40.989043",
    "2021-10-05 13:38:40.989045",
    "2021-10-05 13:38:40.989048",
    "2021-10-05 13:38:40.989050",
    "2021-10-05 13:38:40.989053",
    "2021-10-05 13:38:40.989056",
    "2021-10-05 13:38:40.989058",
    "2021-10-05 13:38:40.989061",
    "2021-10-05 13:38:40.989063",
    "2021-10-05 13:38:40.989066",
    "2021-10-05 13:38:40.989069",
    "2021-10-05 13:38:40.989071",
    "2021-10-05 13:38:40.989073",
    "2021-10-05 13:38:40.989076",
    "2021-10-05 13:38:40.989078",
    "2021-10-05 13:38:40.989081",
    "2021-10-05 13:38:40.989083",
    "2021-10-05 13:38:40.989086",
    "2021-10-05 13:38:40.989089",
    "2021-10-05 13:38:40.989091",
    "2021-10-05 13:38:40.989094",
    "2021-10-05 13:38:40.989096",
    "2021-10-05 13:38:40.989099",
    "2021-10-05 13:38:40.989101",
    "2021-10-05 13:38:40.989104",
    "2021-10-05 13:38:40.989106",
    "2021-10-05 13:38:40.989108",
    "2021-10-05 13:38:40.989111",
    "2021-10-05 13:38:40.989113",
    "2021-10-05 13:38:40.989117",
    "2021-10-05 13:38:40.989120",
    "2021-10-05 13:38:40.989122",
    "2021-10-05 13:38:40.989125",
    "2021-10-05 13:38:40.989127",
    "2021-10-05 13:38:40.989130",
    "2021-10-05 13:38:40.989132",
    "2021-10-05 13:38:40.989135",
    "2021-10-05 13:38:40.989137",
    "2021-10-05 13:38:40.989140",
    "2021-10-05 13:38:40.989142",
    "2021-10-05 13:38:40.989145",
    "2021-10-05 13:38:40.989147",
    "2021-10-05 13:38:40.989150",
    "2021-10-05 13:38:40.989152",
    "2021-10-05 13:38:40.989155",
    "2021-10-05 13:38:40.989157",
    "2021-10-05 13:38:40.989160",
    "2021-10-05 13:38:40.989162",
    "2021-10-05 13:38:40.989165",
    "2021-10-05 13:38:40.989167",
    "2021-10-05 13:38:40.989170",
    "2021-10-05 13:38:40.989172",
    "2021-10-05 13:38:40.989175",
    "2021-10-05 13:38:40.989177",
    "2021-10-05 13:38:40.989180",
    "2021-10-05 13:38:40.989182",
    "2021-10-05 13:38:40.989185",
    "2021-10-05 13:38:40.989188",
    "2021-10-05 13:38:40.989190",
    "2021-10-05 13:38:40.989193",
    "2021-10-05 13:38:40.989195",
    "2021-10-05 13:38:40.989198",
    "2021-10-05 13:38:40.989200",
    "2021-10-05 13:38:40.989203",
    "2021-10-05 13:38:40.989206",
    "2021-10-05 13:38:40.989208",
    "2021-10-05 13:38:40.989211",
    "2021-10-05 13:38:40.989213",
    "2021-10-05 13:38:40.989216",
    "2021-10-05 13:38:40.989218",
    "2021-10-05 13:38:40.989221",
    "2021-10-05 13:38:40.989223",
    "2021-10-05 13:38:40.989226",
    "2021-10-05 13:38:40.989228",
    "2021-10-05 13:38:40.989231",
    "2021-10-05 13:38:40.989233",
    "2021-10-05 13:38:40.989236",
    "2021-10-05 13:38:40.989238",
    "2021-10-05 13:38:40.989240",
    "2021-10-05 13:38:40.989243",
    "2021-10-05 13:38:40.989245",
    "2021-10-05 13:38:40.989248",
    "2021-10-05 13:38:40.989250",
    "2021-10-05 13:38:40.989253",
    "2021-10-05 13:38:40.989255",
    "2021-10-05 13:38:40.989258",
    "2021-10-05 13:38:40.989260",
    "2021-10-05 13:38:40.989263",
    "2021-10-05 13:38:40.989265",
    "2021-10-05 13:38:40.989269",
    "2021-10-05 13:38:40.989272",
    "2021-10-05 13:38:40.989274",
    "2021-10-05 13:38:40.989277",
    "2021-10-05 13:38:40.989279",
    "2021-10-05 13:38:40.989282",
    "2021-10-05 13:38:40.989284",
    "2021-10-05 13:38:40.989287",
    "2021-10-05 13:38:40.989289",
    "2021-10-05 13:38:40.989292",
    "2021-10-05 13:38:40.989294",
    "2021-10-05 13:38:40.989297",
    "2021-10-05 13:38:40.989299",
    "2021-10-05 13:38:40.989302",
    "2021-10-05 13:38:40.989304",
    "2021-10-05 13:38:40.989307",
    "2021-10-05 13:38:40.989310",
    "2021-10-05 13:38:40.989312",
    "2021-10-05 13:38:40.989315",
    "2021-10-05 13:38:40.989317",
    "2021-10-05 13:38:40.989320",
    "2021-10-05 13:38:40.989322",
    "2021-10-05 13:38:40.989325",
    "2021-10-05 13:38:40.989327",
    "2021-10-05 13:38:40.989330",
    "2021-10-05 13:38:40.989332",
    "2021-10-05 13:38:40.989335",
    "2021-10-05 13:38:40.989337",
    "2021-10-05 13:38:40.989339",
    "2021-10-05 13:38:40.989342",
    "2021-10-05 13:38:40.989344",
    "2021-10-05 13:38:40.989347",
    "2021-10-05 13:38:40.989349",
    "2021-10-05 13:38:40.989352",
    "2021-10-05 13:38:40.989354",
    "2021-10-05 13:38:40.989357",
    "2021-10-05 13:38:40.989360",
    "2021-10-05 13:38:40.989362",
    "2021-10-05 13:38:40.989365",
    "2021-10-05 13:38:40.989367",
    "2021-10-05 13:38:40.989370",
    "2021-10-05 13:38:40.989372",
    "2021-10-05 13:38:40.989375",
    "2021-10-05 13:38:40.989377",
    "2021-10-05 13:38:40.989380",
    "2021-10-05 13:38:40.989382",
    "2021-10-05 13:38:40.989385",
    "2021-10-05 13:38:40.989387",
    "2021-10-05 13:38:40.989390",
    "2021-10-05 13:38:40.989392",
    "2021-10-05 13:38:40.989395",
    "2021-10-05 13:38:40.989398",
    "2021-10-05 13:38:40.989400",
    "2021-10-05 13:38:40.989403"
)

for ($i = 0; $i -lt $timeTaken.Count; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $timeTaken[$i]
}

